$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple per-cell D/E updates
$ws.Range("D2").Value = "26.083.40"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.650.45"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.41"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06260"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.75"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.664.80"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.387"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5393"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D17").Value = "26.087.73"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.677"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.34"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.137"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.07"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1211"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.355"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.387"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06010"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.244"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.444"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.399"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.625"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9772"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.379"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.733"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5891"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01595"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8449"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.43"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "1.804.68"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "0.0₈106"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.959"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05215"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.881"
$ws.Range("E51").Value = "  -0.15%  "

# Row 15/16 swap (Litecoin <-> ShibaInu)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₅7960"

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.12"
$ws.Range("E16").Value = "  +1.04%  "

# Row 39/40 swap (FraxShare <-> Maker)
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.075.60"
$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.924"
$ws.Range("E40").Value = "  -2.86%  "
